$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("A2").Formula = "=B1+OneRange"
$ws.Range("A3").Select()
